# Update the cryptos list with freshly scraped data (GitHub Actions style refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $text) {
    # Force the cell to keep a plain-text value even when it looks like a
    # number (e.g. "608.87"), without leaving a lingering text numberformat
    # or quote-prefix style on the cell.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "66.249.36"
$ws.Range("E2").Value = "  +0.48%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.569.16"
$ws.Range("E3").Value = "  +2.64%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
Set-TextCell "D5" "608.87"
$ws.Range("E5").Value = "  +1.23%  "

# Row 6 - Solana
Set-TextCell "D6" "145.38"
$ws.Range("E6").Value = "  +1.74%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.567.45"
$ws.Range("E7").Value = "  +2.62%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.05%  "

# Row 9 - XRP
Set-TextCell "D9" "0.492"
$ws.Range("E9").Value = "  +3.94%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.64%  "

# Row 11 - Toncoin
Set-TextCell "D11" "7.92"
$ws.Range("E11").Value = "  -2.97%  "

# Row 12 - Cardano
Set-TextCell "D12" "0.414"
$ws.Range("E12").Value = "  +0.74%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "4.173.87"
$ws.Range("E13").Value = "  +2.46%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +2.61%  "

# Row 15 - Avalanche
Set-TextCell "D15" "30.02"
$ws.Range("E15").Value = "  -0.85%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.585.76"
$ws.Range("E16").Value = "  +3.08%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "66.332.25"
$ws.Range("E17").Value = "  +0.34%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -0.92%  "

# Row 19 - Uniswap
Set-TextCell "D19" "11.46"
$ws.Range("E19").Value = "  +10.05%  "

# Row 20 - Polkadot
Set-TextCell "D20" "6.23"
$ws.Range("E20").Value = "  +1.12%  "

# Row 21 - Chainlink
Set-TextCell "D21" "14.90"
$ws.Range("E21").Value = "  +1.52%  "

# Row 22 - BitcoinCash
Set-TextCell "D22" "430.17"
$ws.Range("E22").Value = "  +2.60%  "

# Row 23 - Polygon
Set-TextCell "D23" "0.615"
$ws.Range("E23").Value = "  +4.57%  "

# Row 24 - Litecoin
$ws.Range("E24").Value = "  +2.46%  "

# Row 25 - WrappedeETH
$ws.Range("D25").Value = "3.710.91"
$ws.Range("E25").Value = "  +2.61%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.01%  "

# Row 27 - PEPE
$ws.Range("E27").Value = "  +3.70%  "

# Row 28 - PancakeSwap
$ws.Range("E28").Value = "  +2.52%  "

# Row 29 - RenderToken
$ws.Range("E29").Value = "  +0.08%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextCell "D30" "9.11"
$ws.Range("E30").Value = "  -2.21%  "

# Row 31 - Binance-PegBSC-USD
$ws.Range("E31").Value = "  -0.07%  "

# Row 32 & 33 - Fetch.AI and EthereumClassic swap ranking positions
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell "D32" "25.65"
$ws.Range("E32").Value = "  +1.82%  "

$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D33" "1.46"
$ws.Range("E33").Value = "  -1.05%  "

# Row 34 - RenzoRestakedETH
$ws.Range("D34").Value = "3.562.30"
$ws.Range("E34").Value = "  +2.42%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  -6.34%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +2.19%  "

# Row 38 - Aptos
Set-TextCell "D38" "7.89"
$ws.Range("E38").Value = "  +2.71%  "

# Row 39 - NEARProtocol
Set-TextCell "D39" "5.62"
$ws.Range("E39").Value = "  +1.20%  "

# Row 40 - FirstDigitalUSD
$ws.Range("E40").Value = "  -0.07%  "

# Row 41 - Monero
Set-TextCell "D41" "176.29"
$ws.Range("E41").Value = "  +3.91%  "

# Row 42 - Hedera
$ws.Range("E42").Value = "  -1.75%  "

# Row 43 - Filecoin
$ws.Range("E43").Value = "  +2.71%  "

# Row 44 - Mantle
Set-TextCell "D44" "0.898"
$ws.Range("E44").Value = "  +1.00%  "

# Row 45 - Stacks
$ws.Range("E45").Value = "  +1.28%  "

# Row 46 - OKB
Set-TextCell "D46" "46.19"
$ws.Range("E46").Value = "  +2.23%  "

# Row 47 - ONDO
$ws.Range("E47").Value = "  +1.24%  "

# Row 48 - InjectiveProtocol
Set-TextCell "D48" "25.79"
$ws.Range("E48").Value = "  +0.33%  "

# Row 49 - dogwifhat
$ws.Range("E49").Value = "  +2.58%  "

# Row 50 - EnergySwap
Set-TextCell "D50" "23.61"
$ws.Range("E50").Value = "  +9.70%  "

# Row 51 - Cosmos
Set-TextCell "D51" "7.14"
$ws.Range("E51").Value = "  +0.48%  "
